# Update the "Extracurricular Activities" HackerRank blurb:
#  - reword the sentence
#  - apply Verdana/9pt formatting to every run (and the paragraph mark)
#  - split the paragraph in two, moving the (hidden) _GoBack bookmark into
#    a brand-new, empty paragraph inserted right after the reworded text.

$d = $word.ActiveDocument

# Locate the paragraph that still holds the original wording.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Solving code (SQL & Python) challenges on HackerRank*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $rPr = '<w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>'

    $frag = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' + $rPr + '</w:pPr>' +
            '<w:r>' + $rPr + '<w:t xml:space="preserve">Solving Database, SQL, Python challenges on </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r>' + $rPr + '<w:t>HackerRank</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r>' + $rPr + '<w:t>! Just earned the Gold Badge (58/58 challenges solved) for SQL, Silver Badge for Python.</w:t></w:r>' +
            '</w:p>' +
            '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' + $rPr + '</w:pPr>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
            '</w:p>'

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $frag + '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $target.Range.InsertXML($pkg)
}
